# Update "horarios" (schedules) workbook with the latest scrape data.
# This applies the 11:56:32 scrape pass on top of the previous 11:38:09 pass:
#  - header/summary cells are refreshed on the three sheets
#  - a handful of existing rows on the "LP1912" sheet get new values because
#    the scheduling order shuffled once the new timestamps were folded in
#  - ten brand new rows (138-147) are appended to "LP1912" for the extra
#    arrivals captured since the last update
#  - the corresponding rows on the other two sheets pick up the refreshed
#    scrape time + recalculated "Minutos" value

$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, [int]$row, [string]$a, [string]$b, [string]$c, $d, [string]$e)
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 11:56:32"
$ws1.Range("A3").Value = "Total filas: 142"

# Rows whose values swapped/changed position among themselves
Set-Row $ws1 55 "08:14:55" "08:53" "215B_EL PATO" 39 "LP1912"
Set-Row $ws1 56 "08:49:06" "08:53" "23_HERNANDEZ" 4 "LP1912"

Set-Row $ws1 66 "08:49:06" "09:18" "14_ABASTO" 29 "LP1912"
Set-Row $ws1 67 "08:57:42" "09:18" "15X38_ABASTO" 21 "LP1912"

Set-Row $ws1 71 "08:49:06" "09:31" "23_HERNANDEZ" 42 "LP1912"
Set-Row $ws1 72 "08:14:55" "09:31" "16_SANTA ANA" 77 "LP1912"

Set-Row $ws1 90 "09:42:42" "10:29" "14_ABASTO" 47 "LP1912"
Set-Row $ws1 91 "09:42:42" "10:29" "15_ABASTO" 47 "LP1912"

# Rows 117-137 are rewritten (reordered + recalculated Minutos), and
# rows 138-147 are brand new appended rows.
Set-Row $ws1 117 "11:56:32" "11:56" "16_SANTA ANA" 0 "LP1912"
Set-Row $ws1 118 "11:56:32" "11:56" "17_ROMERO" 0 "LP1912"
Set-Row $ws1 119 "11:38:09" "11:58" "17_ROMERO" 20 "LP1912"
Set-Row $ws1 120 "11:01:19" "12:05" "11_ETCHEVERRY" 64 "LP1912"
Set-Row $ws1 121 "11:38:09" "12:06" "11_ETCHEVERRY" 28 "LP1912"
Set-Row $ws1 122 "11:38:09" "12:10" "15_ABASTO" 32 "LP1912"
Set-Row $ws1 123 "11:56:32" "12:10" "16_P MOR-SANTA ANA" 14 "LP1912"
Set-Row $ws1 124 "11:56:32" "12:17" "10_OLMOS" 21 "LP1912"
Set-Row $ws1 125 "11:56:32" "12:22" "215C_EL PATO" 26 "LP1912"
Set-Row $ws1 126 "11:38:09" "12:27" "14_ABASTO" 49 "LP1912"
Set-Row $ws1 127 "11:38:09" "12:31" "23_HERNANDEZ" 53 "LP1912"
Set-Row $ws1 128 "11:56:32" "12:32" "14_ABASTO" 36 "LP1912"
Set-Row $ws1 129 "11:38:09" "12:33" "14_ABASTO" 55 "LP1912"
Set-Row $ws1 130 "11:56:32" "12:34" "15_ABASTO" 38 "LP1912"
Set-Row $ws1 131 "11:56:32" "12:36" "23_HERNANDEZ" 40 "LP1912"
Set-Row $ws1 132 "11:56:32" "12:36" "27_EL RETIRO" 40 "LP1912"
Set-Row $ws1 133 "11:38:09" "12:37" "27_EL RETIRO" 59 "LP1912"
Set-Row $ws1 134 "11:56:32" "12:48" "16_SANTA ANA" 52 "LP1912"
Set-Row $ws1 135 "11:56:32" "12:48" "15X38_ABASTO" 52 "LP1912"
Set-Row $ws1 136 "11:56:32" "13:02" "11_ETCHEVERRY" 66 "LP1912"
Set-Row $ws1 137 "11:38:09" "13:03" "11_ETCHEVERRY" 85 "LP1912"
Set-Row $ws1 138 "11:56:32" "13:04" "215C_EL PATO" 68 "LP1912"
Set-Row $ws1 139 "11:38:09" "13:13" "16_SANTA ANA" 95 "LP1912"
Set-Row $ws1 140 "11:56:32" "13:17" "10_OLMOS" 81 "LP1912"
Set-Row $ws1 141 "11:56:32" "13:24" "23_HERNANDEZ" 88 "LP1912"
Set-Row $ws1 142 "11:56:32" "13:25" "16_P MOR-SANTA ANA" 89 "LP1912"
Set-Row $ws1 143 "11:56:32" "13:32" "14_ABASTO" 96 "LP1912"
Set-Row $ws1 144 "11:56:32" "13:33" "215A_EL PATO" 97 "LP1912"
Set-Row $ws1 145 "11:56:32" "13:47" "225_GOMEZ" 111 "LP1912"
Set-Row $ws1 146 "11:56:32" "13:52" "11_ETCHEVERRY" 116 "LP1912"
Set-Row $ws1 147 "11:56:32" "13:54" "15_ABASTO" 118 "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 11:56:32"

$ws2.Cells.Item(26, 1).Value = "11:56:32"
$ws2.Cells.Item(26, 4).Value = 26

$ws2.Cells.Item(27, 1).Value = "11:56:32"
$ws2.Cells.Item(27, 4).Value = 68

$ws2.Cells.Item(28, 1).Value = "11:56:32"
$ws2.Cells.Item(28, 4).Value = 97

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 11:56:32"

$ws3.Cells.Item(26, 1).Value = "11:56:32"
$ws3.Cells.Item(26, 4).Value = 76

$ws3.Cells.Item(27, 1).Value = "11:56:32"
$ws3.Cells.Item(27, 4).Value = 85
